$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: Provincia 6 -> 4, Localizacion -> "El CAR de la liga" ---
$ws.Range("F6").Value = 4.0
$ws.Range("J6").Value = "El CAR de la liga"

# --- Row 132: Provincia 1 -> 7, Localizacion -> "Ceneguita" ---
$ws.Range("F132").Value = 7.0
$ws.Range("J132").Value = "Ceneguita"

# --- Row 133: full record replaced ---
# Fecha (date-like text) needs to be forced as text so it doesn't get
# auto-converted into a date serial number; reset the style afterwards so no
# explicit cell style sticks around (matches the original, style-less cell).
$ws.Range("A133").NumberFormat = "@"
$ws.Range("A133").Value = "09/10/2021"
$ws.Range("A133").Style = "Normal"
$ws.Range("B133").Value = "21:29:00"
$ws.Range("C133").Value = 100.0
$ws.Range("D133").Value = 4.3
$ws.Range("E133").Value = "Deformación Interna"
$ws.Range("F133").Value = 5.0
$ws.Range("G133").Value = 100.0
$ws.Range("H133").Value = -500.4
$ws.Range("I133").Value = 1.0
$ws.Range("J133").Value = "Marcial Fallas"

# --- Row 134: full record replaced ---
$ws.Range("A134").NumberFormat = "@"
$ws.Range("A134").Value = "09/10/2021"
$ws.Range("A134").Style = "Normal"
$ws.Range("B134").Value = "21:42:00"
$ws.Range("C134").Value = 100.0
$ws.Range("D134").Value = 3.2
$ws.Range("E134").Value = "Deformación Interna"
$ws.Range("F134").Value = 1.0
$ws.Range("G134").Value = 100.4
$ws.Range("H134").Value = -903.5
$ws.Range("I134").Value = 1.0
$ws.Range("J134").Value = "Maxi Pali"

# --- Row 135: brand new record appended at the bottom ---
# Strip styles first so the new row's cells don't inherit the column
# default styles (the target row has no explicit per-cell style, same as
# rows 133/134 above it).
$ws.Range("A135:J135").Style = "Normal"
$ws.Range("A135").NumberFormat = "@"
$ws.Range("A135").Value = "09/10/2021"
$ws.Range("A135").Style = "Normal"
$ws.Range("B135").Value = "09:55:12"
$ws.Range("C135").Value = 33.0
$ws.Range("D135").Value = 5.6
$ws.Range("E135").Value = "Intra placa"
$ws.Range("F135").Value = 3.0
$ws.Range("G135").Value = -312.0
$ws.Range("H135").Value = -4444.0
$ws.Range("I135").Value = 2.0
$ws.Range("J135").Value = "Lomas de Ayarco"
